$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R11").Value = "LU"
$ws.Range("S20").Value = "RL"
$ws.Range("N17").Value = "WU"
$ws.Range("D18").Value = "SU"
$ws.Range("C8").Value = "MR"
$ws.Range("C2").Value = "BD"
$ws.Range("G5").Value = "CL"
$ws.Range("N4").Value = "GD"
$ws.Range("S7").Value = "EL"

$ws.Range("S7").Select()
